# Apply the TMDB testdata changes:
#  - rename Sheet1 -> favorites and fill it with favorites data
#  - add a new sheet "watchlist" with the same rows but a "watchlist" header

$wb = $excel.ActiveWorkbook

# --- Rename existing Sheet1 to "favorites" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "favorites"

$headers = @("media_type", "media_id", "favorite")
$data = @(
    @("movie", 496243, $true),
    @("movie", 129,    $true),
    @("movie", 278,    $true),
    @("movie", 155,    $true),
    @("movie", 27205,  $true),
    @("tv",    1396,   $true),
    @("tv",    66732,  $true),
    @("tv",    1399,   $true),
    @("tv",    2316,   $true),
    @("tv",    1429,   $true)
)

for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws1.Cells.Item(1, $c + 1).Value = $headers[$c]
}

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws1.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

$ws1.Range("B2").Select()

# --- Add a new "watchlist" sheet right after "favorites", with the same data
#     but a different header for column C ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "watchlist"

$headers2 = @("media_type", "media_id", "watchlist")
for ($c = 0; $c -lt $headers2.Length; $c++) {
    $ws2.Cells.Item(1, $c + 1).Value = $headers2[$c]
}

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

$ws2.Range("A2:C11").Select()

$ws1.Activate()
